$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price / 1h volume data.
# For Price (column D) values that are plain numeric-looking strings,
# a leading apostrophe is used so Excel stores them as text (matching
# the original inlineStr cell type) instead of auto-converting to a number.

$ws.Range("D2").Value = '58.407.77'
$ws.Range("E2").Value = '  -1.31%  '

$ws.Range("D3").Value = '2.591.28'
$ws.Range("E3").Value = '  -1.97%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '''518.80'
$ws.Range("E5").Value = '  -1.69%  '

$ws.Range("D6").Value = '''141.73'
$ws.Range("E6").Value = '  -2.16%  '

$ws.Range("D8").Value = '''0.564'
$ws.Range("E8").Value = '  -0.85%  '

$ws.Range("D9").Value = '2.613.74'
$ws.Range("E9").Value = '  -1.67%  '

$ws.Range("D10").Value = '''6.49'
$ws.Range("E10").Value = '  -2.53%  '

$ws.Range("E11").Value = '  -2.81%  '

$ws.Range("D12").Value = '''0.332'
$ws.Range("E12").Value = '  -1.39%  '

$ws.Range("E13").Value = '  +0.45%  '

$ws.Range("D14").Value = '3.040.22'
$ws.Range("E14").Value = '  -2.18%  '

$ws.Range("D15").Value = '58.400.48'

$ws.Range("D16").Value = '''20.33'
$ws.Range("E16").Value = '  -3.29%  '

$ws.Range("E17").Value = '  -2.76%  '

$ws.Range("D18").Value = '2.594.15'
$ws.Range("E18").Value = '  -2.71%  '

$ws.Range("D19").Value = '''338.63'
$ws.Range("E19").Value = '  -1.02%  '

$ws.Range("D20").Value = '''4.32'
$ws.Range("E20").Value = '  -3.09%  '

$ws.Range("D21").Value = '''10.22'
$ws.Range("E21").Value = '  -3.28%  '

$ws.Range("D22").Value = '''6.43'
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").Value = '''65.42'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").Value = '''0.168'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").Value = '''0.403'
$ws.Range("E26").Value = '  -3.89%  '

$ws.Range("D27").Value = '''0.996'
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").Value = '2.698.97'
$ws.Range("E28").Value = '  -2.02%  '

$ws.Range("D29").Value = '''7.04'
$ws.Range("E29").Value = '  -3.02%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0739'
$ws.Range("E30").Value = '  -7.55%  '

$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").Value = '''0.998'
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").Value = '''6.08'
$ws.Range("E32").Value = '  -6.12%  '

$ws.Range("E33").Value = '  -2.53%  '

$ws.Range("D34").Value = '''18.77'
$ws.Range("E34").Value = '  -1.09%  '

$ws.Range("D35").Value = '''149.39'
$ws.Range("E35").Value = '  -0.51%  '

$ws.Range("D36").Value = '''3.98'
$ws.Range("E36").Value = '  -6.20%  '

$ws.Range("E37").Value = '  -6.07%  '

$ws.Range("D38").Value = '''0.857'
$ws.Range("E38").Value = '  -2.46%  '

$ws.Range("D39").Value = '''36.26'
$ws.Range("E39").Value = '  -0.86%  '

$ws.Range("E40").Value = '  -1.86%  '

$ws.Range("D41").Value = '''0.832'
$ws.Range("E41").Value = '  -10.53%  '

$ws.Range("D42").Value = '''3.52'
$ws.Range("E42").Value = '  -3.81%  '

$ws.Range("D43").Value = '''0.997'
$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("D44").Value = '''274.17'
$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("D45").Value = '''0.601'
$ws.Range("E45").Value = '  -0.52%  '

$ws.Range("D46").Value = '''10.68'
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").Value = '''0.0948'
$ws.Range("E47").Value = '  -2.70%  '

$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '''0.0520'
$ws.Range("E48").Value = '  -3.64%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''18.63'
$ws.Range("E49").Value = '  -4.25%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''4.65'
$ws.Range("E50").Value = '  -3.31%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.971.34'
$ws.Range("E51").Value = '  -4.01%  '
